$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.411.21"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.847.49"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.45"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07606"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2931"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.50"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07743"
$ws.Range("D12").Value = "1.848.50"
$ws.Range("E12").Value = "  -6.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.004"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001084"
$ws.Range("E14").Value = "  +8.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6796"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.79"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "2.100.05"
$ws.Range("E17").Value = "  -7.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.174"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "29.424.33"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.75"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.45"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.466"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.44"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1396"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.357"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.64"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.463"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.301"
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05582"
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.103"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.031"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.842"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.156"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7098"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.586"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "1.232.98"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01801"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.772"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.427"
$ws.Range("E41").Value = "  +5.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9060"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.82"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.08"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("E46").Value = "  +3.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.201"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4021"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.956"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("E51").Value = "  -0.52%  "
